$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17 updates
$ws.Range("J17").Value = "29.01.25 05:59"
$ws.Range("N17").Value = "strieme: engl. weal > wheelless > radlos > ratlos: ratlos im Stall (dt. Idiom?)"

# Row 18 updates
$ws.Range("H18").Value = "wendig"
$ws.Range("J18").Value = "29.01.25 06:14"
$ws.Range("N18").Value = "engl.: astute > Stute"
$ws.Range("O18").Value = 6
$ws.Range("P18").Value = 0.05624824224242992
